$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new columns "I0" (I) and "IF" (J) after the existing "IP" (H) column.
# Copy H1's formatting (bold, bordered, centered header style) onto the two new
# header cells, then overwrite their text so the style is reused rather than a
# brand-new style being minted.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Per-row data for columns I ("I0") and J ("IF"), rows 2-73 (one entry per row,
# in order).
$iVals = @(9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,6,8,8,9,8,8,8,8,8,7,8,8,8,9,8,8,8,8,8,8,8,8,8,8,9,7,8,9,7,7,7,7,8,8,7,9,7,5,10,7,8,7,7,9,9,7,8,4,5,9,7,9,7,5,4)
$jVals = @(9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,6,8,8,9,8,8,8,8,8,7,8,8,8,9,8,8,8,8,8,8,8,8,8,8,9,8,8,9,7,7,7,7,8,8,7,9,8,6,10,7,8,8,7,9,9,7,8,4,5,9,7,9,7,5,4)

for ($r = 2; $r -le 73; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iVals[$idx]
    $ws.Cells.Item($r, 10).Value = $jVals[$idx]
}
